# ykodama/ddbj_packages MIGS.ba.microbial.4.0.xlsx
# "culture_collection" column (AH, MIxS) is removed again per INSDC2017 review.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire "culture_collection" column (AH). This removes its header cell,
# drops the now-unreferenced shared string, and shifts every later column one slot left.
$ws.Range("AH1").EntireColumn.Delete() | Out-Null

# Column delete does not re-anchor the per-cell header comments on its own, so walk the
# row of field-description comments and pull each one in from the column that used to
# sit one to the right (AI..CS before the delete -> AH..CR after it).
$ws.Range("AH15").Comment.Text('concentration of diether lipids; can include multiple types of diether lipids') | Out-Null
$ws.Range("AI15").Comment.Text('concentration of dissolved carbon dioxide') | Out-Null
$ws.Range("AJ15").Comment.Text('concentration of dissolved hydrogen') | Out-Null
$ws.Range("AK15").Comment.Text('dissolved inorganic carbon concentration') | Out-Null
$ws.Range("AL15").Comment.Text('concentration of dissolved organic carbon') | Out-Null
$ws.Range("AM15").Comment.Text('dissolved organic nitrogen concentration measured as; total dissolved nitrogen - NH4 - NO3 - NO2') | Out-Null
$ws.Range("AN15").Comment.Text('concentration of dissolved oxygen') | Out-Null
$ws.Range("AO15").Comment.Text('Traits like antibiotic resistance/xenobiotic degration phenotypes/converting phage genes') | Out-Null
$ws.Range("AP15").Comment.Text('Estimated size of genome') | Out-Null
$ws.Range("AQ15").Comment.Text('Plasmids that have significance phenotypic consequence') | Out-Null
$ws.Range("AR15").Comment.Text('measurement of glucosidase activity') | Out-Null
$ws.Range("AS15").Comment.Text('Health or disease status of sample at time of collection') | Out-Null
$ws.Range("AT15").Comment.Text('The natural (as opposed to laboratory) host to the organism from which the sample was obtained. Use the full taxonomic name, eg, "Homo sapiens".') | Out-Null
$ws.Range("AU15").Comment.Text('NCBI taxonomy ID of the host, e.g. 9606') | Out-Null
$ws.Range("AV15").Comment.Text('Describes the physical, environmental and/or local geographical source of the biological sample from which the sample was derived.') | Out-Null
$ws.Range("AW15").Comment.Text('A locus tag prefix required for an annotated genome, http://www.ddbj.nig.ac.jp/sub/locus_tag-e.html') | Out-Null
$ws.Range("AX15").Comment.Text('concentration of magnesium') | Out-Null
$ws.Range("AY15").Comment.Text('measurement of mean friction velocity') | Out-Null
$ws.Range("AZ15").Comment.Text('measurement of mean peak friction velocity') | Out-Null
$ws.Range("BA15").Comment.Text('methane (gas) amount or concentration at the time of sampling') | Out-Null
$ws.Range("BB15").Comment.Text('any other measurement performed or parameter collected, that is not listed here') | Out-Null
$ws.Range("BC15").Comment.Text('concentration of n-alkanes; can include multiple n-alkanes') | Out-Null
$ws.Range("BD15").Comment.Text('concentration of nitrate') | Out-Null
$ws.Range("BE15").Comment.Text('concentration of nitrite') | Out-Null
$ws.Range("BF15").Comment.Text('concentration of nitrogen (total)') | Out-Null
$ws.Range("BG15").Comment.Text('concentration of organic carbon') | Out-Null
$ws.Range("BH15").Comment.Text('concentration of organic matter') | Out-Null
$ws.Range("BI15").Comment.Text('concentration of organic nitrogen') | Out-Null
$ws.Range("BJ15").Comment.Text('total count of any organism per gram or volume of sample, should include name of organism followed by count; can include multiple organism counts') | Out-Null
$ws.Range("BK15").Comment.Text('oxygenation status of sample') | Out-Null
$ws.Range("BL15").Comment.Text('concentration of particulate organic carbon') | Out-Null
$ws.Range("BM15").Comment.Text('To what is the entity pathogenic') | Out-Null
$ws.Range("BN15").Comment.Text('type of perturbation, e.g. chemical administration, physical disturbance, etc., coupled with time that perturbation occurred; can include multiple perturbation types') | Out-Null
$ws.Range("BO15").Comment.Text('concentration of petroleum hydrocarbon') | Out-Null
$ws.Range("BP15").Comment.Text('pH measurement') | Out-Null
$ws.Range("BQ15").Comment.Text('concentration of phaeopigments; can include multiple phaeopigments') | Out-Null
$ws.Range("BR15").Comment.Text('concentration of phosphate') | Out-Null
$ws.Range("BS15").Comment.Text('concentration of phospholipid fatty acids; can include multiple values') | Out-Null
$ws.Range("BT15").Comment.Text('concentration of potassium') | Out-Null
$ws.Range("BU15").Comment.Text('pressure to which the sample is subject, in atmospheres') | Out-Null
$ws.Range("BV15").Comment.Text('redox potential, measured relative to a hydrogen cell, indicating oxidation or reduction potential') | Out-Null
$ws.Range("BW15").Comment.Text('Aerobic or anaerobic') | Out-Null
$ws.Range("BX15").Comment.Text('salinity measurement') | Out-Null
$ws.Range("BY15").Comment.Text('Method or device employed for collecting sample') | Out-Null
$ws.Range("BZ15").Comment.Text('Processing applied to the sample during or after isolation') | Out-Null
$ws.Range("CA15").Comment.Text('Amount or size of sample (volume, mass or area) that was collected') | Out-Null
$ws.Range("CB15").Comment.Text('duration for which sample was stored') | Out-Null
$ws.Range("CC15").Comment.Text('location at which sample was stored, usually name of a specific freezer/room') | Out-Null
$ws.Range("CD15").Comment.Text('temperature at which sample was stored, e.g. -80') | Out-Null
$ws.Range("CE15").Comment.Text('volume (mL) or weight (g) of sample processed for DNA extraction') | Out-Null
$ws.Range("CF15").Comment.Text('concentration of silicate') | Out-Null
$ws.Range("CG15").Comment.Text('sodium concentration') | Out-Null
$ws.Range("CH15").Comment.Text('unique identifier assigned to a material sample used for extracting nucleic acids, and subsequent sequencing. The identifier can refer either to the original material collected or to any derived sub-samples.') | Out-Null
$ws.Range("CI15").Comment.Text('Information about the genetic distinctness of the lineage (eg., biovar, serovar)') | Out-Null
$ws.Range("CJ15").Comment.Text('concentration of sulfate') | Out-Null
$ws.Range("CK15").Comment.Text('concentration of sulfide') | Out-Null
$ws.Range("CL15").Comment.Text('temperature of the sample at time of sampling') | Out-Null
$ws.Range("CM15").Comment.Text('total carbon content') | Out-Null
$ws.Range("CN15").Comment.Text('total nitrogen content of the sample') | Out-Null
$ws.Range("CO15").Comment.Text('Definition for soil: total organic C content of the soil units of g C/kg soil. Definition otherwise: total organic carbon content') | Out-Null
$ws.Range("CP15").Comment.Text('Feeding position in food chain (eg., chemolithotroph)') | Out-Null
$ws.Range("CQ15").Comment.Text('turbidity measurement') | Out-Null
$ws.Range("CR15").Comment.Text('water content measurement') | Out-Null

# The old last column (CS) has nothing left to shift into it; drop its orphaned comment.
$ws.Range("CS15").Comment.Delete() | Out-Null
